$wb = $excel.ActiveWorkbook

# --- Rename "Movies" sheet to "Movie_list" (Admin Edit Movie Feature) ---
$wsMovies = $wb.Worksheets.Item("Movies")
$wsMovies.Name = "Movie_list"

# --- Add two new movie rows to the Movie_list sheet ---
# Columns: A=Title B=Genre C=Length D=Cast E=Director F=Admin Rating
#          G=Language H=Number of Shows in a day I=First show
#          J=Interval Time K=Gap Between Shows L=Timings M=Capacity

$wsMovies.Range("B2").Value = "Fantasy"
$wsMovies.Range("C2").Value = "2hr 30m"
$wsMovies.Range("D2").Value = "Robert Jr."
$wsMovies.Range("E2").Value = "Stan lee"
$wsMovies.Range("F2").Value = 3.5
$wsMovies.Range("G2").Value = "Eng"
$wsMovies.Range("H2").Value = 4
$wsMovies.Range("I2").Value = "8h 0m"
$wsMovies.Range("J2").Value = "0h 30m"
$wsMovies.Range("K2").Value = "0h 15min"
$wsMovies.Range("L2").Value = "1-2 2-3"
$wsMovies.Range("M2").Value = 2

$wsMovies.Range("B3").Value = "Fantasy"
$wsMovies.Range("C3").Value = "2hr 30m"
$wsMovies.Range("D3").Value = "Robert Jr."
$wsMovies.Range("E3").Value = "Stan lee"
$wsMovies.Range("F3").Value = 3.5
$wsMovies.Range("G3").Value = "Eng"
$wsMovies.Range("H3").Value = 4
$wsMovies.Range("I3").Value = "8h 0m"
$wsMovies.Range("J3").Value = "0h 30m"
$wsMovies.Range("K3").Value = "0h 15min"
$wsMovies.Range("L3").Value = "1-2 2-3"
$wsMovies.Range("M3").Value = 2

$wsMovies.Range("A2").Value = "Test1"
$wsMovies.Range("A3").Value = "Test3"

# --- Update the selection on the UserDetails sheet ---
$wsUserDetails = $wb.Worksheets.Item("UserDetails")
$wsUserDetails.Activate()
$wsUserDetails.Range("F9").Select()

# --- Re-activate Movie_list and update its selection (keeps it the active tab) ---
$wsMovies.Activate()
$wsMovies.Range("G8").Select()
